# Update "想去人数" (interest count) values in column F across all four sheets.
# Values scraped per-sheet: row number -> new F-column value.
$wb = $excel.ActiveWorkbook

$updates = @{
    4 = 493
    5 = 2129
    7 = 7973
    8 = 273
    9 = 59
    11 = 236
    12 = 1767
    13 = 1544
    14 = 1315
    15 = 192
    16 = 4046
    17 = 703
    18 = 32
    19 = 1117
    22 = 6230
    23 = 351
    25 = 4231
    26 = 705
    27 = 1947
    29 = 305
    30 = 1035
    34 = 51
    35 = 320
    36 = 1154
    38 = 1874
    39 = 111
    40 = 411
    41 = 155
    42 = 1145
    45 = 45
}
$ws = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$updates = @{
    11 = 679
    12 = 376
    19 = 343
    20 = 174
    22 = 90
    30 = 273
}
$ws = $wb.Worksheets.Item("演出")
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$updates = @{
    6 = 1560
    8 = 3102
    9 = 943
    10 = 1085
    11 = 1276
    12 = 1588
}
$ws = $wb.Worksheets.Item("本地生活")
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$updates = @{
    3 = 1560
    4 = 493
    6 = 3102
    7 = 7973
    8 = 273
    9 = 59
    10 = 943
    13 = 236
    14 = 1767
    15 = 1544
    16 = 1315
    17 = 679
    18 = 192
    19 = 1588
    20 = 4046
    21 = 376
    23 = 703
    24 = 32
    25 = 1117
    28 = 6230
    29 = 351
    30 = 705
    31 = 1947
    33 = 305
    36 = 51
    37 = 90
    38 = 320
    39 = 1154
    41 = 1874
    42 = 111
    43 = 411
    44 = 1145
    47 = 273
}
$ws = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
